$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save" - copy formatting from the neighboring header (G1, style s="1")
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New data column H2:H9 with the "Save" flag values (plain numbers, no special style)
$saveValues = @(1, 0, 0, 1, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
